$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "2024-02-09 Friday" "2024-02-10 Saturday"
Replace-Text "77÷4=" "50÷3="
Replace-Text "80÷3=" "53÷2="
Replace-Text "59÷4=" "26÷4="
Replace-Text "62÷3=" "56÷4="
Replace-Text "56÷2=" "11÷9="
Replace-Text "58÷3=" "84÷2="
Replace-Text "83÷2=" "65÷3="
Replace-Text "45÷2=" "92÷9="
Replace-Text "70÷6=" "28÷2="
Replace-Text "82÷6=" "91÷3="
Replace-Text "21÷4=" "86÷3="
Replace-Text "61÷8=" "53÷9="
Replace-Text "59÷7=" "55÷4="
Replace-Text "27÷5=" "44÷9="
Replace-Text "51÷4=" "39÷9="
Replace-Text "39÷5=" "87÷9="
Replace-Text "93÷6=" "43÷3="
Replace-Text "14÷5=" "98÷6="
Replace-Text "20÷8=" "29÷4="
Replace-Text "88÷8=" "88÷7="
Replace-Text "71÷6=" "50÷7="
Replace-Text "69÷4=" "44÷3="
Replace-Text "76÷2=" "60÷8="
Replace-Text "77÷7=" "71÷9="
Replace-Text "56÷8=" "40÷8="
